# StorageCAPEX.xlsx edit:
#  - Rename the CAPEX header from "CAPEX [EUR/kWh]" to "CAPEX [EUR/kW]"
#  - Move the active selection from B4 to B2
#  - Restore/minimize the workbook window position & size

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text for CAPEX column (kWh -> kW)
$ws.Range("B1").Value = "CAPEX [EUR/kW]"

# Move the selected cell to B2
$ws.Range("B2").Select()

# Reposition / minimize the workbook window to match the saved view state
$win = $excel.ActiveWindow
$win.Top = 2580
$win.Left = 2580
$win.Width = 14400
$win.Height = 7270
$win.WindowState = -4140
